$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","F","G","H","K","L","M","N","O")

$data = @(
    @(0.4640172698507286,0.3352467479746366,0.03713129615863409,0.7237954068354213,0.5654879152002366,0.7046062645174516,0.2344323194235045,0.2975902557459023,0.1715587280925135,1.699601298836583,2.515919276992392),
    @(0.4282631982735552,0.3369559749467186,0.03436536845873661,0.7240507523788651,0.5672943804798578,0.7086917066676364,0.2044361969667108,0.2941627243645826,0.1641598477432105,1.715745469639874,2.528071116980612),
    @(0.4064267507648367,0.3380675140652194,0.03265123221821398,0.7245971059525687,0.5687296117110066,0.7114618991301072,0.1859594257016397,0.2922071844856404,0.1596888094830078,1.726173716295994,2.536763336162494),
    @(0.3975582094535355,0.3385361169604391,0.03194875477886683,0.7249178163121641,0.5693964975257657,0.7126566702474193,0.1784156134481378,0.2914478408541612,0.1578850460230754,1.730553021586636,2.540615248542366),
    @(0.3960874250831523,0.338614873848087,0.03183187110766283,0.7249769961370944,0.569512187894297,0.7128590434892956,0.1771621144824707,0.291324023368233,0.1575866370908123,1.731288037363338,2.541273571391457),
    @(0.4063070242403342,0.3380737704284265,0.03264177432512838,0.7246010339547624,0.5687382734431381,0.7114777453109724,0.1858577448104199,0.2921967915304648,0.1596644093366351,1.726232251888829,2.536814029874819),
    @(0.4516654204919348,0.33582322824725,0.03618091117813549,0.7238026565841125,0.5660431279150657,0.7059606568783465,0.2241021789450599,0.2963775765413459,0.1689927519031933,1.705060810937481,2.519853886614342),
    @(0.5415152596551422,0.3319008218616162,0.04299432657980873,0.7253243031719805,0.5633443465369297,0.6972146396492178,0.2986139881502368,0.3057548215564054,0.1878510194241088,1.667632680648981,2.496353262859884),
    @(0.6080526127312851,0.3293160946131266,0.04792178198705699,0.7283204950564368,0.5629381486211713,0.6920480096366646,0.3530440763183265,0.3133594902833323,0.2020455460567163,1.642622643838122,2.485026478564293),
    @(0.6384311163348286,0.3282042732030277,0.05014619301184098,0.7300906436314136,0.5630957240756516,0.6899700013653103,0.3777341155911529,0.3169736965693204,0.2085756467486846,1.631784337438654,2.481161349380983),
    @(0.6499499565371423,0.3277924212502725,0.05098603113340516,0.7308194157070105,0.5632046111054052,0.6892221930554143,0.3870730428748459,0.3183644838023838,0.2110587913086093,1.627757604991698,2.479882674976295),
    @(0.6474685051822178,0.3278807135464348,0.05080526866389334,0.7306598633079773,0.5631789716418467,0.6893815096139662,0.385062219941716,0.3180639685769222,0.2105235443668221,1.628621386934718,2.480149836724081),
    @(0.6393784784588377,0.3281702062718246,0.05021533727834537,0.7301494291647614,0.5631036959847791,0.6899076957409562,0.3785026512749425,0.3170876737830497,0.2087797302431227,1.631451502488726,2.48105244606171),
    @(0.6344250568642735,0.3283487223302082,0.04985366084231657,0.7298443832821278,0.5630639965040416,0.6902350882985502,0.3744833230680911,0.3164925488410688,0.2077129354411298,1.633195122909882,2.481629403263867),
    @(0.6060694659284991,0.3293900396975893,0.04777606395862222,0.7282129968261941,0.5629347386213226,0.6921892860498673,0.3514290615183882,0.313126400359053,0.2016202434092307,1.643341797023853,2.485304966579378),
    @(0.5887019857795508,0.3300452204770856,0.04649711606352724,0.7273164119194391,0.56294311927455,0.6934578198324317,0.3372675973400305,0.3111009590216298,0.1979011480845543,1.649704479262693,2.487889438447525),
    @(0.5787230973348017,0.3304280873913648,0.04575989052661811,0.7268390500953359,0.5629801647350803,0.6942130829448274,0.3291156948053526,0.3099505553605866,0.1957688972947764,1.653414847569652,2.489497152553326),
    @(0.5753462294406972,0.3305587551450309,0.04551000349029266,0.7266840096867782,0.5629982428001483,0.6944732072264941,0.3263544824496307,0.3095635548805262,0.195048139119578,1.654679825831345,2.490062317383206),
    @(0.590549709659399,0.3299748520959049,0.04663342899679179,0.7274078885751507,0.5629388920579146,0.6933201295251337,0.3387757965147387,0.3113150627698928,0.1982963417443955,1.649021911311709,2.487601775099051),
    @(0.6417543095919598,0.3280849265887653,0.05038868251985207,0.7302977703474411,0.5631244706594316,0.6897520818363603,0.3804296489093417,0.3173738344168413,0.2092916513837082,1.630618125573459,2.480782309270097),
    @(0.6753075130396269,0.3269031910262221,0.05282837978433008,0.7325271523060977,0.5635326221809009,0.6876479604629822,0.4075905038442613,0.3214627464742108,0.2165379040909556,1.619041916203168,2.47740343756729),
    @(0.6573917146165229,0.3275290256525771,0.05152761367135383,0.7313061487607087,0.563288540631163,0.6887501488291861,0.3931001208585485,0.3192686310042347,0.212664986113495,1.625179020627035,2.479108221132662),
    @(0.5897143352391936,0.3300066463427846,0.04657180796370142,0.7273664133043098,0.5629407027964817,0.6933822984254903,0.3380939716749367,0.3112182226959135,0.1981176562811271,1.649330337129916,2.48773144808095),
    @(0.5171145277089977,0.3329096110966887,0.04116480182418059,0.7245825388630465,0.5637975594416673,0.6993592268542415,0.2785103765435508,0.3030921844243295,0.1826893768902451,1.677320684818519,2.501667127458333)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $ws.Range("$($cols[$j])$row").Value = $data[$i][$j]
    }
}
